# netCrypto.xlsx update:
#  - T2 (USD Amount for the BTC deposit row) changes from 100000 to 10
#  - Row 3 (the second "Wiretransfer" deposit entry) is deleted entirely,
#    which also drops the now-unused "Wiretransfer" shared string and
#    shrinks the used range from A1:AB3 to A1:AB2
#  - Selection ends up on T3 (what used to be just below the deleted row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount on row 2
$ws.Range("T2").Value = 10

# Remove row 3 completely (E3, N3, P3, T3 and the "Wiretransfer" strings go with it)
$ws.Rows("3:3").Delete()

# Leave the selection where row 3 used to start
$ws.Range("T3").Select()
